# Applies the "Fixed issue for no leafs provided in spreadsheet" edit.
#
# Summary of content changes:
#  - "General Configuration Details" sheet: update the example Fabric
#    Name / Configlet Prefix / Fabric Identifier values to the new
#    generic "avd" / "AVD-Demo" sample values, and add explanatory
#    Notes for the DNS/NTP server rows.
#  - "Tenants" sheet: row 3's "Virtual Address Type" value is renamed
#    from "Virtual Router Address" to "IP Virtual Router Address".

$wb = $excel.ActiveWorkbook

# ---- General Configuration Details ----------------------------------
$general = $wb.Worksheets.Item("General Configuration Details")

$general.Range("B2").Value = "AVD-Demo"
$general.Range("C2").Value = "#Root container"

$general.Range("B3").Value = "avd"
$general.Range("C3").Value = "#prefx for configlet i.e. If device name is 'lf1' the configlet name will be 'avd_lf1'"

$general.Range("B4").Value = "avd"
$general.Range("C4").Value = "#Device filter"

$general.Range("C7").Value = "DNS Servers separated by a comma"
$general.Range("C8").Value = "NTP Servers separated by a comma"

# ---- Tenants ----------------------------------------------------------
$tenants = $wb.Worksheets.Item("Tenants")
$tenants.Range("I3").Value = "IP Virtual Router Address"

# Re-select the cells Excel left active/focused after the edit session.
$general.Activate()
$general.Range("C18").Select()
